$wb = $excel.ActiveWorkbook

# --- Update the conversion text on sheet "Hoja1" ---
$hoja1 = $wb.Worksheets.Item("Hoja1")

$oldText = $hoja1.Range("A1").Value2
$newText = $oldText.Replace(
    "1000 Bs = 7.27 = 29587.13 pesos",
    "1000 Bs = 7.0 = 28455.4 pesos"
).Replace(
    "29587.13 pesos = 7.24 = 945.92 Bs",
    "28455.4 pesos = 6.98 = 956.72 Bs"
)
$hoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$tasas = $wb.Worksheets.Item("tasas")

$tasas.Range("N10").Value = 142.82
$tasas.Range("O10").Value = 4064
$tasas.Range("N12").Value = 4077.85
$tasas.Range("O12").Value = 137.105
